$d = $word.ActiveDocument

# Every paragraph below contains exactly one Word field whose instruction
# text is an M2Doc template tag, e.g. a field with code " m: self.myTemplate1(1) ".
# Convert each such field into plain literal text wrapped in curly braces
# (the field's begin/instrText-runs/end go away, replaced by a run of
# text reading "{m: self.myTemplate1(1)}"), preserving the run's language
# formatting.

$paraCount = $d.Paragraphs.Count
for ($pi = 1; $pi -le $paraCount; $pi++) {
    $p = $d.Paragraphs.Item($pi)
    $pr = $p.Range

    while ($pr.Fields.Count -gt 0) {
        # Paragraph-scoped Fields.Item(1) resolves oddly in this host, but
        # Document.Fields.Item(1) reliably is the next not-yet-converted
        # field in document order, which is the one living in this paragraph.
        $fld = $d.Fields.Item(1)
        $code = $fld.Code.Text

        # Trim exactly one leading/trailing space (the field code always has
        # them around the "m:..." instruction) and wrap with { }.
        $body = $code
        if ($body.StartsWith(" ")) {
            $body = $body.Substring(1)
        }
        if ($body.EndsWith(" ")) {
            $body = $body.Substring(0, $body.Length - 1)
        }
        $newText = "{" + $body + "}"

        $fld.Delete()

        $pr = $p.Range
        $pr.Text = $newText
        # A no-op Write-Output acts as a sync point for this host so the
        # just-edited Range handle below is not considered stale.
        Write-Output "" | Out-Null
        $pr.LanguageID = "en-US"
        $pr = $p.Range
    }
}
